# Trade #3 closed at 2026-02-17 13:33:23 - unknown UNKNOWN +0.000%
#
# Updates the Summary, Strategy Status, All Trades, and MarketMaking
# sheets to reflect a newly closed trade (#3) for the MarketMaking
# strategy.

$wb = $excel.ActiveWorkbook

# ---- Summary sheet ----
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1199.91   # Current Capital
$wsSummary.Range("B4").Value = -0.09     # Total P&L $
$wsSummary.Range("B5").Value = -0.6      # Total P&L %
$wsSummary.Range("B6").Value = 3         # Total Trades
$wsSummary.Range("B8").Value = 3         # Losing Trades

# ---- Strategy Status sheet (MarketMaking row) ----
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 99.91      # Capital
$wsStatus.Range("D4").Value = 3          # Trades
$wsStatus.Range("E4").Value = -0.09      # P&L $
$wsStatus.Range("F4").Value = -0.09      # P&L %

# ---- Append new trade row (#3) to "All Trades" and "MarketMaking" sheets ----
foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $newRow = 4

    $ws.Cells.Item($newRow, 1).Value = 3                 # A: Trade #

    # Date/time columns look like dates/times to the auto-detector, so force
    # them to stay plain text (matches the other rows' stored type).
    $ws.Cells.Item($newRow, 2).NumberFormat = "@"
    $ws.Cells.Item($newRow, 2).Value = "2026-02-17"       # B: Date
    $ws.Cells.Item($newRow, 3).NumberFormat = "@"
    $ws.Cells.Item($newRow, 3).Value = "13:33:17"         # C: Time

    $ws.Cells.Item($newRow, 4).Value = "MarketMaking"     # D: Strategy
    $ws.Cells.Item($newRow, 5).Value = "UP"                # E: Side
    $ws.Cells.Item($newRow, 6).Value = 0.15                # F: Entry Price
    $ws.Cells.Item($newRow, 7).Value = 0.09                # G: Exit Price
    $ws.Cells.Item($newRow, 8).Value = "CLOSED"            # H: Status
    $ws.Cells.Item($newRow, 9).Value = -40                 # I: P&L %
    $ws.Cells.Item($newRow, 10).Value = -0.06              # J: P&L $
    $ws.Cells.Item($newRow, 11).Value = 99.91              # K: Capital After
    $ws.Cells.Item($newRow, 12).Value = 0                  # L: Entry Slippage (bps)
    $ws.Cells.Item($newRow, 13).Value = 0                  # M: Exit Slippage (bps)
    $ws.Cells.Item($newRow, 14).Value = 0.6                # N: Confidence
    $ws.Cells.Item($newRow, 15).Value = "Normal spread capture: 19600 bps"  # O: Entry Reason
    $ws.Cells.Item($newRow, 16).Value = "early_exit"       # P: Exit Reason
    $ws.Cells.Item($newRow, 17).Value = 0.13                # Q: Duration (min)
}
